# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy the formatting used by the other header cells (e.g. G1)
# so the new header shares the same style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Row values for the new "Save" column (0/1 flags).
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 0
    37 = 0
}

foreach ($row in 2..37) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
